# Updates cryptos list values (price / 1h volume change) per the latest scrape.
# Price-column (D) numeric-looking strings are written with a leading apostrophe
# so Excel keeps them as literal text (preserving formats like "69.006.61",
# "0.0450", "0.999", thousand-dot grouping, etc.) instead of coercing them to
# numbers; the style is then reset to Normal so no stray number-format/quote-prefix
# styling is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'69.006.61"
$ws.Range("E2").Value = '  +3.16%  '
# Row 3: Ethereum
$ws.Range("D3").Value = "'3.721.97"
$ws.Range("E3").Value = '  +2.18%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.05%  '
# Row 5: BNB
$ws.Range("D5").Value = "'614.63"
$ws.Range("E5").Value = '  +9.78%  '
# Row 6: Solana
$ws.Range("D6").Value = "'192.95"
$ws.Range("E6").Value = '  +13.95%  '
# Row 7: XRP
$ws.Range("E7").Value = '  +4.02%  '
# Row 8: USDC
$ws.Range("E8").Value = '  +0.12%  '
# Row 9: Cardano
$ws.Range("E9").Value = '  +5.01%  '
# Row 10: Dogecoin
$ws.Range("E10").Value = '  +2.31%  '
# Row 11: Avalanche
$ws.Range("D11").Value = "'60.42"
$ws.Range("E11").Value = '  +21.76%  '
# Row 12: ShibaInu
$ws.Range("E12").Value = '  +2.26%  '
# Row 13: Polkadot
$ws.Range("D13").Value = "'10.49"
$ws.Range("E13").Value = '  +2.11%  '
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'4.315.32"
$ws.Range("E14").Value = '  +2.25%  '
# Row 15: WrappedEther
$ws.Range("D15").Value = "'3.722.00"
$ws.Range("E15").Value = '  +1.80%  '
# Row 16: Chainlink
$ws.Range("E16").Value = '  +2.92%  '
# Row 17: Polygon
$ws.Range("E17").Value = '  +4.79%  '
# Row 18: TRON
$ws.Range("E18").Value = '  +1.28%  '
# Row 19: Uniswap
$ws.Range("D19").Value = "'12.98"
$ws.Range("E19").Value = '  +3.14%  '
# Row 20: WrappedBTC
$ws.Range("D20").Value = "'68.915.68"
$ws.Range("E20").Value = '  +3.16%  '
# Row 21: BitcoinCash
$ws.Range("D21").Value = "'413.09"
$ws.Range("E21").Value = '  +3.78%  '
# Row 22: PancakeSwap
$ws.Range("E22").Value = '  +5.79%  '
# Row 23: Litecoin
$ws.Range("D23").Value = "'90.25"
$ws.Range("E23").Value = '  +4.79%  '
# Row 24: ImmutableX
$ws.Range("D24").Value = "'3.11"
$ws.Range("E24").Value = '  +4.60%  '
# Row 25: InternetComputer(DFINITY)->RenderToken
$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").Value = "'11.44"
$ws.Range("E25").Value = '  +9.85%  '
# Row 26: RenderToken->InternetComputer(DFINITY)
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = "'13.12"
$ws.Range("E26").Value = '  +5.39%  '
# Row 27: Toncoin
$ws.Range("E27").Value = '  +3.57%  '
# Row 28: LEO
$ws.Range("D28").Value = "'6.05"
$ws.Range("E28").Value = '  +1.25%  '
# Row 29: Filecoin
$ws.Range("D29").Value = "'9.75"
$ws.Range("E29").Value = '  +5.83%  '
# Row 30: EthereumClassic
$ws.Range("D30").Value = "'32.97"
$ws.Range("E30").Value = '  +2.96%  '
# Row 31: NEARProtocol
$ws.Range("D31").Value = "'7.83"
$ws.Range("E31").Value = '  +5.55%  '
# Row 32: Cosmos
$ws.Range("E32").Value = '  +4.49%  '
# Row 33: Hedera
$ws.Range("E33").Value = '  +8.49%  '
# Row 34: InjectiveProtocol
$ws.Range("D34").Value = "'46.37"
$ws.Range("E34").Value = '  +10.31%  '
# Row 35: Bittensor
$ws.Range("D35").Value = "'638.87"
$ws.Range("E35").Value = '  +10.35%  '
# Row 36: OKB
$ws.Range("D36").Value = "'67.59"
$ws.Range("E36").Value = '  +5.27%  '
# Row 37: PEPE->TheGraph
$ws.Range("B37").Value = 'TheGraph'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D37").Value = "'0.418"
$ws.Range("E37").Value = '  +8.47%  '
# Row 38: TheGraph->PEPE
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = "'0.0₃0835"
$ws.Range("E38").Value = '  -4.79%  '
# Row 39: Dai
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = '  +0.00%  '
# Row 40: FirstDigitalUSD
$ws.Range("E40").Value = '  +0.10%  '
# Row 41: Kaspa
$ws.Range("E41").Value = '  +8.42%  '
# Row 42: ThetaToken
$ws.Range("E42").Value = '  +5.52%  '
# Row 43: VeChain
$ws.Range("D43").Value = "'0.0450"
$ws.Range("E43").Value = '  +5.08%  '
# Row 44: Fetch.AI
$ws.Range("E44").Value = '  +5.86%  '
# Row 45: Maker
$ws.Range("D45").Value = "'2.933.14"
$ws.Range("E45").Value = '  +9.67%  '
# Row 46: Stellar
$ws.Range("E46").Value = '  +6.43%  '
# Row 47: THORChain
$ws.Range("D47").Value = "'9.32"
$ws.Range("E47").Value = '  +4.28%  '
# Row 48: WEMIXToken
$ws.Range("D48").Value = "'2.73"
$ws.Range("E48").Value = '  +3.59%  '
# Row 49: Monero
$ws.Range("D49").Value = "'144.24"
$ws.Range("E49").Value = '  +2.31%  '
# Row 50: ApeXProtocol
$ws.Range("D50").Value = "'3.10"
$ws.Range("E50").Value = '  -1.01%  '
# Row 51: dogwifhat
$ws.Range("E51").Value = '  -11.52%  '

# Clear the quote-prefix styling picked up by the text-forced Price cells above
# so the cells end up style-equivalent to the rest of the untouched column.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
